# Auto-generated edit script applying the Gilgamesh_Profits workbook diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1720
$ws.Range("I6").Value = 1720
$ws.Range("K6").Value = 5160
$ws.Range("M6").Value = -5048
$ws.Range("H12").Value = 179.25
$ws.Range("I12").Value = 195
$ws.Range("J12").Value = 69
$ws.Range("K12").Value = 195
$ws.Range("L12").Value = 69
$ws.Range("M12").Value = -25
$ws.Range("N12").Value = -409
$ws.Range("H19").Value = 683.2727
$ws.Range("I19").Value = 497.2
$ws.Range("J19").Value = 838.3333
$ws.Range("K19").Value = 497.2
$ws.Range("L19").Value = 838.3333
$ws.Range("M19").Value = -322.2
$ws.Range("N19").Value = -1188.3333
$ws.Range("H29").Value = 5464.6665
$ws.Range("J29").Value = 5464.6665
$ws.Range("L29").Value = 16393.9995
$ws.Range("N29").Value = -16955.9995
$ws.Range("H38").Value = 5252.737
$ws.Range("I38").Value = 1821
$ws.Range("J38").Value = 11135.714
$ws.Range("K38").Value = 5463
$ws.Range("L38").Value = 33407.142
$ws.Range("M38").Value = -5091
$ws.Range("N38").Value = -34151.142
$ws.Range("H40").Value = 4034.3572
$ws.Range("J40").Value = 4874.125
$ws.Range("L40").Value = 4874.125
$ws.Range("N40").Value = -5224.125
$ws.Range("H42").Value = 206.125
$ws.Range("I42").Value = 80
$ws.Range("J42").Value = 281.8
$ws.Range("K42").Value = 240
$ws.Range("L42").Value = 845.4000000000001
$ws.Range("M42").Value = -10
$ws.Range("N42").Value = -1305.4
$ws.Range("H70").Value = 787.5
$ws.Range("I70").Value = 685
$ws.Range("K70").Value = 2055
$ws.Range("M70").Value = -1785
$ws.Range("H73").Value = 787.5
$ws.Range("I73").Value = 685
$ws.Range("K73").Value = 2055
$ws.Range("M73").Value = -1119
$ws.Range("H87").Value = 157249.75
$ws.Range("J87").Value = 264499.5
$ws.Range("L87").Value = 264499.5
$ws.Range("N87").Value = -266995.5
$ws.Range("H90").Value = 157249.75
$ws.Range("J90").Value = 264499.5
$ws.Range("L90").Value = 793498.5
$ws.Range("N90").Value = -805978.5
$ws.Range("H94").Value = 851.75
$ws.Range("I94").Value = 851.75
$ws.Range("K94").Value = 851.75
$ws.Range("M94").Value = -400.75
$ws.Range("H98").Value = 1922.4062
$ws.Range("I98").Value = 1972.1072
$ws.Range("K98").Value = 1972.1072
$ws.Range("M98").Value = -474.1071999999999
$ws.Range("H113").Value = 5053.75
$ws.Range("I113").Value = 1802.5
$ws.Range("K113").Value = 1802.5
$ws.Range("M113").Value = 1451.5
$ws.Range("H116").Value = 11570.714
$ws.Range("J116").Value = 13749
$ws.Range("L116").Value = 13749
$ws.Range("N116").Value = -20633
$ws.Range("H122").Value = 1922.4062
$ws.Range("I122").Value = 1972.1072
$ws.Range("K122").Value = 5916.321599999999
$ws.Range("M122").Value = -3466.321599999999
$ws.Range("H125").Value = 772.8
$ws.Range("I125").Value = 772.8
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 6955.2
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -4495.2
$ws.Range("N125").ClearContents()
$ws.Range("H129").Value = 2507.5
$ws.Range("I129").Value = 1612
$ws.Range("J129").Value = 4000
$ws.Range("K129").Value = 4836
$ws.Range("L129").Value = 12000
$ws.Range("M129").Value = 164
$ws.Range("N129").Value = -22000
$ws.Range("H132").Value = 4765.795
$ws.Range("I132").Value = 4274.1943
$ws.Range("K132").Value = 12822.5829
$ws.Range("M132").Value = -10292.5829
$ws.Range("H137").Value = 2587.8
$ws.Range("I137").Value = 2399.818
$ws.Range("J137").Value = 2817.5557
$ws.Range("K137").Value = 7199.454000000001
$ws.Range("L137").Value = 8452.667099999999
$ws.Range("M137").Value = -4649.454000000001
$ws.Range("N137").Value = -13552.6671
$ws.Range("H138").Value = 4848.909
$ws.Range("J138").Value = 4968.8965
$ws.Range("L138").Value = 14906.6895
$ws.Range("N138").Value = -25186.6895
$ws.Range("H141").Value = 3221.75
$ws.Range("I141").Value = 3221.75
$ws.Range("K141").Value = 9665.25
$ws.Range("M141").Value = -4485.25

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 60
$ws.Range("I5").Value = 60
$ws.Range("K5").Value = 60
$ws.Range("M5").Value = 52
$ws.Range("H32").Value = 3666.625
$ws.Range("I32").Value = 3704.7737
$ws.Range("K32").Value = 3704.7737
$ws.Range("M32").Value = -3417.7737
$ws.Range("H61").Value = 11556.75
$ws.Range("I61").Value = 8606.5
$ws.Range("K61").Value = 8606.5
$ws.Range("M61").Value = -8394.5
$ws.Range("H74").Value = 35721124
$ws.Range("I74").Value = 4623.0835
$ws.Range("J74").Value = 62508504
$ws.Range("K74").Value = 4623.0835
$ws.Range("L74").Value = 62508504
$ws.Range("M74").Value = -3749.0835
$ws.Range("N74").Value = -62510252
$ws.Range("H77").Value = 35721124
$ws.Range("I77").Value = 4623.0835
$ws.Range("J77").Value = 62508504
$ws.Range("K77").Value = 23115.4175
$ws.Range("L77").Value = 312542520
$ws.Range("M77").Value = -18747.4175
$ws.Range("N77").Value = -312551256
$ws.Range("H88").Value = 7217.25
$ws.Range("I88").Value = 1899.8334
$ws.Range("J88").Value = 10407.7
$ws.Range("K88").Value = 1899.8334
$ws.Range("L88").Value = 10407.7
$ws.Range("M88").Value = -1493.8334
$ws.Range("N88").Value = -11219.7
$ws.Range("H91").Value = 7217.25
$ws.Range("I91").Value = 1899.8334
$ws.Range("J91").Value = 10407.7
$ws.Range("K91").Value = 1899.8334
$ws.Range("L91").Value = 10407.7
$ws.Range("M91").Value = -495.8334
$ws.Range("N91").Value = -13215.7
$ws.Range("H102").Value = 2160.9697
$ws.Range("I102").Value = 1925.4642
$ws.Range("K102").Value = 1925.4642
$ws.Range("M102").Value = -303.4641999999999
$ws.Range("H110").Value = 1617
$ws.Range("I110").Value = 1571.4348
$ws.Range("K110").Value = 1571.4348
$ws.Range("M110").Value = 473.5652
$ws.Range("H122").Value = 2826.7273
$ws.Range("I122").Value = 2268.6667
$ws.Range("J122").Value = 3803.3333
$ws.Range("K122").Value = 6806.000100000001
$ws.Range("L122").Value = 11409.9999
$ws.Range("M122").Value = -4356.000100000001
$ws.Range("N122").Value = -16309.9999
$ws.Range("H132").Value = 3074.48
$ws.Range("J132").Value = 3012.3
$ws.Range("L132").Value = 9036.900000000001
$ws.Range("N132").Value = -14096.9
$ws.Range("H136").Value = 11556.75
$ws.Range("I136").Value = 8606.5
$ws.Range("K136").Value = 25819.5
$ws.Range("M136").Value = -23269.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 60
$ws.Range("I4").Value = 60
$ws.Range("K4").Value = 60
$ws.Range("M4").Value = 55
$ws.Range("H20").Value = 27781938
$ws.Range("I20").Value = 32055810
$ws.Range("K20").Value = 32055810
$ws.Range("M20").Value = -32055563
$ws.Range("H80").Value = 288
$ws.Range("I80").Value = 347
$ws.Range("J80").Value = 273.25
$ws.Range("K80").Value = 347
$ws.Range("L80").Value = 273.25
$ws.Range("M80").Value = 651
$ws.Range("N80").Value = -2269.25
$ws.Range("H83").Value = 288
$ws.Range("I83").Value = 347
$ws.Range("J83").Value = 273.25
$ws.Range("K83").Value = 1735
$ws.Range("L83").Value = 1366.25
$ws.Range("M83").Value = 3257
$ws.Range("N83").Value = -11350.25
$ws.Range("H87").Value = 75000
$ws.Range("J87").Value = 75000
$ws.Range("L87").Value = 75000
$ws.Range("N87").Value = -77496
$ws.Range("H90").Value = 75000
$ws.Range("J90").Value = 75000
$ws.Range("L90").Value = 225000
$ws.Range("N90").Value = -237480
$ws.Range("H94").Value = 38096670
$ws.Range("I94").Value = 38096670
$ws.Range("K94").Value = 38096670
$ws.Range("M94").Value = -38096219
$ws.Range("H99").Value = 338002
$ws.Range("I99").Value = 1000010
$ws.Range("K99").Value = 1000010
$ws.Range("M99").Value = -998512
$ws.Range("H105").Value = 10834978
$ws.Range("I105").Value = 589234.1
$ws.Range("K105").Value = 589234.1
$ws.Range("M105").Value = -587487.1
$ws.Range("H107").Value = 2405264.5
$ws.Range("I107").Value = 2565431.8
$ws.Range("K107").Value = 2565431.8
$ws.Range("M107").Value = -2563511.8
$ws.Range("H134").Value = 2116.3713
$ws.Range("I134").Value = 1520.9546
$ws.Range("J134").Value = 3124
$ws.Range("K134").Value = 4562.8638
$ws.Range("L134").Value = 9372
$ws.Range("M134").Value = -2027.8638
$ws.Range("N134").Value = -14442

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 33334450
$ws.Range("I16").Value = 41667780
$ws.Range("K16").Value = 41667780
$ws.Range("M16").Value = -41667493
$ws.Range("H31").Value = 2282.319
$ws.Range("I31").Value = 1742.1765
$ws.Range("J31").Value = 2588.4
$ws.Range("K31").Value = 1742.1765
$ws.Range("L31").Value = 2588.4
$ws.Range("M31").Value = -1447.1765
$ws.Range("N31").Value = -3178.4
$ws.Range("H34").Value = 2282.319
$ws.Range("I34").Value = 1742.1765
$ws.Range("J34").Value = 2588.4
$ws.Range("K34").Value = 1742.1765
$ws.Range("L34").Value = 2588.4
$ws.Range("M34").Value = -1540.1765
$ws.Range("N34").Value = -2992.4
$ws.Range("H58").Value = 3814.25
$ws.Range("I58").Value = 3166.6667
$ws.Range("J58").Value = 4202.8
$ws.Range("K58").Value = 3166.6667
$ws.Range("L58").Value = 4202.8
$ws.Range("M58").Value = -2963.6667
$ws.Range("N58").Value = -4608.8
$ws.Range("H74").Value = 77509
$ws.Range("J74").Value = 77509
$ws.Range("L74").Value = 77509
$ws.Range("N74").Value = -79257
$ws.Range("H77").Value = 77509
$ws.Range("J77").Value = 77509
$ws.Range("L77").Value = 232527
$ws.Range("N77").Value = -241263
$ws.Range("H86").Value = 5043.222
$ws.Range("I86").Value = 3174.75
$ws.Range("J86").Value = 6538
$ws.Range("K86").Value = 3174.75
$ws.Range("L86").Value = 6538
$ws.Range("M86").Value = -2051.75
$ws.Range("N86").Value = -8784
$ws.Range("H89").Value = 5043.222
$ws.Range("I89").Value = 3174.75
$ws.Range("J89").Value = 6538
$ws.Range("K89").Value = 15873.75
$ws.Range("L89").Value = 32690
$ws.Range("M89").Value = -10257.75
$ws.Range("N89").Value = -43922
$ws.Range("H113").Value = 33334450
$ws.Range("I113").Value = 41667780
$ws.Range("K113").Value = 41667780
$ws.Range("M113").Value = -41665610
$ws.Range("H121").Value = 49999
$ws.Range("J121").Value = 49999
$ws.Range("L121").Value = 49999
$ws.Range("N121").Value = -52619
$ws.Range("H122").Value = 3139.75
$ws.Range("I122").Value = 3227.2
$ws.Range("K122").Value = 9681.599999999999
$ws.Range("M122").Value = -7231.599999999999
$ws.Range("H132").Value = 3198.9048
$ws.Range("I132").Value = 2709.7058
$ws.Range("J132").Value = 5278
$ws.Range("K132").Value = 8129.117400000001
$ws.Range("L132").Value = 15834
$ws.Range("M132").Value = -5599.117400000001
$ws.Range("N132").Value = -20894
$ws.Range("H134").Value = 5661.3
$ws.Range("I134").Value = 4013.2942
$ws.Range("J134").Value = 15000
$ws.Range("K134").Value = 12039.8826
$ws.Range("L134").Value = 45000
$ws.Range("M134").Value = -9504.882599999999
$ws.Range("N134").Value = -50070
$ws.Range("H136").Value = 3814.25
$ws.Range("I136").Value = 3166.6667
$ws.Range("J136").Value = 4202.8
$ws.Range("K136").Value = 9500.000100000001
$ws.Range("L136").Value = 12608.4
$ws.Range("M136").Value = -6950.000100000001
$ws.Range("N136").Value = -17708.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 750.36
$ws.Range("J5").Value = 901.3570999999999
$ws.Range("L5").Value = 2704.0713
$ws.Range("N5").Value = -2928.0713
$ws.Range("H62").Value = 9499.666999999999
$ws.Range("J62").Value = 9499.666999999999
$ws.Range("L62").Value = 28499.001
$ws.Range("N62").Value = -29871.001
$ws.Range("H65").Value = 9499.666999999999
$ws.Range("J65").Value = 9499.666999999999
$ws.Range("L65").Value = 85497.003
$ws.Range("N65").Value = -92361.003
$ws.Range("H68").Value = 1473740.9
$ws.Range("J68").Value = 1669907
$ws.Range("L68").Value = 5009721
$ws.Range("N68").Value = -5011343
$ws.Range("H71").Value = 1473740.9
$ws.Range("J71").Value = 1669907
$ws.Range("L71").Value = 15029163
$ws.Range("N71").Value = -15037275
$ws.Range("H80").Value = 2559
$ws.Range("I80").Value = 2189
$ws.Range("K80").Value = 6567
$ws.Range("M80").Value = -5631
$ws.Range("H83").Value = 2559
$ws.Range("I83").Value = 2189
$ws.Range("K83").Value = 19701
$ws.Range("M83").Value = -15021
$ws.Range("H107").Value = 5320.857
$ws.Range("J107").Value = 8311
$ws.Range("L107").Value = 24933
$ws.Range("N107").Value = -28773
$ws.Range("H127").Value = 902.8333
$ws.Range("J127").Value = 902.8333
$ws.Range("L127").Value = 2708.4999
$ws.Range("N127").Value = -12628.4999
$ws.Range("H132").Value = 5481.472
$ws.Range("I132").Value = 3525.2354
$ws.Range("J132").Value = 7231.7896
$ws.Range("K132").Value = 31727.1186
$ws.Range("L132").Value = 65086.1064
$ws.Range("M132").Value = -29197.1186
$ws.Range("N132").Value = -70146.1064
$ws.Range("H135").Value = 750.36
$ws.Range("J135").Value = 901.3570999999999
$ws.Range("L135").Value = 8112.2139
$ws.Range("N135").Value = -13182.2139
$ws.Range("H139").Value = 7456.4346
$ws.Range("I139").Value = 8432.5
$ws.Range("K139").Value = 25297.5
$ws.Range("M139").Value = -20157.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H51").Value = 38999
$ws.Range("J51").Value = 38999
$ws.Range("L51").Value = 38999
$ws.Range("N51").Value = -40017
$ws.Range("H97").Value = 1782.5555
$ws.Range("I97").Value = 1784.6666
$ws.Range("K97").Value = 1784.6666
$ws.Range("M97").Value = -1288.6666
$ws.Range("H107").Value = 3136.4167
$ws.Range("J107").Value = 3874.2222
$ws.Range("L107").Value = 3874.2222
$ws.Range("N107").Value = -7714.2222
$ws.Range("H113").Value = 83336900
$ws.Range("I113").Value = 115388020
$ws.Range("K113").Value = 115388020
$ws.Range("M113").Value = -115385850
$ws.Range("H122").Value = 60611396
$ws.Range("I122").Value = 47623140
$ws.Range("J122").Value = 83340840
$ws.Range("K122").Value = 142869420
$ws.Range("L122").Value = 250022520
$ws.Range("M122").Value = -142866970
$ws.Range("N122").Value = -250027420
$ws.Range("H126").Value = 35715960
$ws.Range("I126").Value = 83334664
$ws.Range("J126").Value = 1937.5
$ws.Range("K126").Value = 250003992
$ws.Range("L126").Value = 5812.5
$ws.Range("M126").Value = -250001522
$ws.Range("N126").Value = -10752.5
$ws.Range("H132").Value = 1863.7273
$ws.Range("I132").Value = 1837.9231
$ws.Range("J132").Value = 1901
$ws.Range("K132").Value = 5513.7693
$ws.Range("L132").Value = 5703
$ws.Range("M132").Value = -2983.7693
$ws.Range("N132").Value = -10763

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3250
$ws.Range("J7").Value = 3000
$ws.Range("L7").Value = 3000
$ws.Range("N7").Value = -3224
$ws.Range("H16").Value = 3458.0833
$ws.Range("I16").Value = 3156.8572
$ws.Range("K16").Value = 3156.8572
$ws.Range("M16").Value = -2986.8572
$ws.Range("H22").Value = 3053.7144
$ws.Range("I22").Value = 1492.6666
$ws.Range("J22").Value = 4224.5
$ws.Range("K22").Value = 1492.6666
$ws.Range("L22").Value = 4224.5
$ws.Range("M22").Value = -1197.6666
$ws.Range("N22").Value = -4814.5
$ws.Range("H27").Value = 3053.7144
$ws.Range("I27").Value = 1492.6666
$ws.Range("J27").Value = 4224.5
$ws.Range("K27").Value = 1492.6666
$ws.Range("L27").Value = 4224.5
$ws.Range("M27").Value = -1385.6666
$ws.Range("N27").Value = -4438.5
$ws.Range("H68").Value = 2022.8
$ws.Range("I68").Value = 2264.182
$ws.Range("K68").Value = 2264.182
$ws.Range("M68").Value = -1515.182
$ws.Range("H71").Value = 2022.8
$ws.Range("I71").Value = 2264.182
$ws.Range("K71").Value = 11320.91
$ws.Range("M71").Value = -7576.91
$ws.Range("H100").Value = 2310.3333
$ws.Range("I100").Value = 3332
$ws.Range("K100").Value = 3332
$ws.Range("M100").Value = -2791
$ws.Range("H122").Value = 6352.4
$ws.Range("I122").Value = 5000
$ws.Range("J122").Value = 7254
$ws.Range("K122").Value = 15000
$ws.Range("L122").Value = 21762
$ws.Range("M122").Value = -12550
$ws.Range("N122").Value = -26662
$ws.Range("H126").Value = 3250
$ws.Range("J126").Value = 3000
$ws.Range("L126").Value = 9000
$ws.Range("N126").Value = -13940
$ws.Range("H132").Value = 7748.2666
$ws.Range("I132").Value = 7909.9287
$ws.Range("J132").Value = 7482
$ws.Range("K132").Value = 23729.7861
$ws.Range("L132").Value = 22446
$ws.Range("M132").Value = -21199.7861
$ws.Range("N132").Value = -27506
$ws.Range("H136").Value = 6424
$ws.Range("I136").Value = 4671.0835
$ws.Range("K136").Value = 14013.2505
$ws.Range("M136").Value = -11463.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 29850
$ws.Range("J54").Value = 29750.666
$ws.Range("L54").Value = 29750.666
$ws.Range("N54").Value = -30790.666
$ws.Range("H62").Value = 10749
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H64").Value = 827691.4
$ws.Range("J64").Value = 827691.4
$ws.Range("L64").Value = 827691.4
$ws.Range("N64").Value = -828187.4
$ws.Range("H65").Value = 10749
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H67").Value = 827691.4
$ws.Range("J67").Value = 827691.4
$ws.Range("L67").Value = 827691.4
$ws.Range("N67").Value = -829407.4
$ws.Range("H70").Value = 39127.285
$ws.Range("J70").Value = 38998.2
$ws.Range("L70").Value = 38998.2
$ws.Range("N70").Value = -39628.2
$ws.Range("H73").Value = 39127.285
$ws.Range("J73").Value = 38998.2
$ws.Range("L73").Value = 38998.2
$ws.Range("N73").Value = -41182.2
$ws.Range("H100").Value = 76925200
$ws.Range("J100").Value = 125002200
$ws.Range("L100").Value = 250004400
$ws.Range("N100").Value = -250005482
$ws.Range("H122").Value = 7355648
$ws.Range("I122").Value = 2796.7693
$ws.Range("J122").Value = 31252414
$ws.Range("K122").Value = 8390.3079
$ws.Range("L122").Value = 93757242
$ws.Range("M122").Value = -5940.3079
$ws.Range("N122").Value = -93762142
$ws.Range("H132").Value = 3968.077
$ws.Range("I132").Value = 3870.7188
$ws.Range("J132").Value = 4413.143
$ws.Range("K132").Value = 11612.1564
$ws.Range("L132").Value = 13239.429
$ws.Range("M132").Value = -9082.1564
$ws.Range("N132").Value = -18299.429
$ws.Range("H136").Value = 47636360
$ws.Range("I136").Value = 55573530
$ws.Range("J136").Value = 13350
$ws.Range("K136").Value = 166720590
$ws.Range("L136").Value = 40050
$ws.Range("M136").Value = -166718040
$ws.Range("N136").Value = -45150
$ws.Range("H137").Value = 119996.336
$ws.Range("J137").Value = 119996.336
$ws.Range("L137").Value = 119996.336
$ws.Range("N137").Value = -130196.336
